$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.979.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.380.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.959.02"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.396.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.065.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.79%  "
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("E26").Value = "  +3.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "165.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.415.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0765"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.16%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.458.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.27%  "
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("E51").Value = "  -2.00%  "

